{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Replace the four \"Key Achievements\" bullet paragraphs (and their Heading3\n// label) that previously described software-development work with the\n// research-leadership / policy-impact framing from the commit.\n\nconst replacements = [\n  {\n    find: \"Software Development and Innovation\",\n    replace: \"Research Leadership and Policy Impact\"\n  },\n  {\n    find: \"\\u2022 Conceived and deployed redistricting software used by thousands of analysts nationwide\",\n    replace: \"\\u2022 Regular expert testimony and consultation on research methodology for journalists, elected officials, and community leaders\"\n  },\n  {\n    find: \"\\u2022 Developed boundary estimation system using incomplete data without ML requirements\",\n    replace: \"\\u2022 Research analysis used in court cases addressing housing, redistricting, and community development with rigorous methodology\"\n  },\n  {\n    find: \"\\u2022 Created econometric simulation platform for humanitarian intervention modeling\",\n    replace: \"\\u2022 Conceived and deployed cloud-based analytical software used by thousands of researchers nationwide for community-focused research\"\n  },\n  {\n    find: \"\\u2022 Built comprehensive survey operations platform from RFP through deployment\",\n    replace: \"\\u2022 Developed research frameworks and methodologies that became industry standards for community development and policy analysis\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# The document is already open as $word.ActiveDocument.\n#\n# Replace the four \"Key Achievements\" bullet paragraphs (and their Heading3\n# label) that previously described software-development work with the\n# research-leadership / policy-impact framing from the commit.\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute(\n        $findText,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $replaceText,# ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-DocText \"Software Development and Innovation\" \"Research Leadership and Policy Impact\"\nReplace-DocText \"Conceived and deployed redistricting software used by thousands of analysts nationwide\" \"Regular expert testimony and consultation on research methodology for journalists, elected officials, and community leaders\"\nReplace-DocText \"Developed boundary estimation system using incomplete data without ML requirements\" \"Research analysis used in court cases addressing housing, redistricting, and community development with rigorous methodology\"\nReplace-DocText \"Created econometric simulation platform for humanitarian intervention modeling\" \"Conceived and deployed cloud-based analytical software used by thousands of researchers nationwide for community-focused research\"\nReplace-DocText \"Built comprehensive survey operations platform from RFP through deployment\" \"Developed research frameworks and methodologies that became industry standards for community development and policy analysis\"\n"}
